$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.440.55'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '2.160.94'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.64'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.22'
$ws.Range("E7").Value = '  +4.13%  '
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0858'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.03'
$ws.Range("E12").Value = '  +4.44%  '
$ws.Range("D13").Value = '2.481.55'
$ws.Range("E13").Value = '  +3.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.20'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.814'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.55'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").Value = '2.159.00'
$ws.Range("E17").Value = '  +3.65%  '
$ws.Range("D18").Value = '39.387.41'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.79'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.11'
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("D21").Value = '0.0₃0851'
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.40'
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +5.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.44'
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.86'
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("E31").Value = '  +7.08%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.61'
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.76'
$ws.Range("E35").Value = '  +9.03%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.40'
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '103.89'
$ws.Range("E40").Value = '  +2.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0230'
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.82'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("D43").Value = '1.539.60'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.83'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0925'
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.26'
$ws.Range("E47").Value = '  +3.49%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.10'
$ws.Range("E48").Value = '  +5.65%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.76'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").Value = '2.364.81'
$ws.Range("E50").Value = '  +3.21%  '
$ws.Range("E51").Value = '  -0.11%  '
